$wb = $excel.ActiveWorkbook

# --- Sheet: Matriz_Resultados (sign matrix: 1 win / -1 loss / 0 tie) ---
$ws1 = $wb.Worksheets.Item("Matriz_Resultados")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 0
$ws1.Range("G2").Value = 0
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 0
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 0
$ws1.Range("E4").Value = 0
$ws1.Range("I4").Value = 0
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("B7").Value = 0
$ws1.Range("D9").Value = 0

# --- Sheet: P_valores (DM test p-values) ---
$ws2 = $wb.Worksheets.Item("P_valores")
$ws2.Range("C2").Value = 0.02242931537559456
$ws2.Range("D2").Value = 0.00783345273120517
$ws2.Range("E2").Value = 0.003518479473069824
$ws2.Range("F2").Value = [double]"3.636635927950493E-07"
$ws2.Range("G2").Value = 0.003703522545972948
$ws2.Range("H2").Value = 0.01820238026848409
$ws2.Range("I2").Value = 0.02441417733121787
$ws2.Range("J2").Value = [double]"2.661584286300922E-10"
$ws2.Range("B3").Value = 0.02242931537559456
$ws2.Range("D3").Value = 0.005191625890059726
$ws2.Range("E3").Value = 0.002468448345302354
$ws2.Range("F3").Value = [double]"6.31289682306857E-07"
$ws2.Range("G3").Value = 0.0001670144785241856
$ws2.Range("H3").Value = 0.0005751431850828936
$ws2.Range("I3").Value = 0.06538983234925433
$ws2.Range("J3").Value = [double]"6.841640587396114E-10"
$ws2.Range("B4").Value = 0.00783345273120517
$ws2.Range("C4").Value = 0.005191625890059726
$ws2.Range("E4").Value = 0.01518547225210143
$ws2.Range("F4").Value = 0.0130729465972379
$ws2.Range("G4").Value = 0.1707566629146371
$ws2.Range("H4").Value = 0.5848326098805563
$ws2.Range("I4").Value = 0.004547430365912764
$ws2.Range("J4").Value = [double]"3.257989211746803E-09"
$ws2.Range("B5").Value = 0.003518479473069824
$ws2.Range("C5").Value = 0.002468448345302354
$ws2.Range("D5").Value = 0.01518547225210143
$ws2.Range("F5").Value = 0.5840467020102458
$ws2.Range("G5").Value = 0.8728760011712222
$ws2.Range("H5").Value = 0.4503888688957607
$ws2.Range("I5").Value = 0.001294854920740995
$ws2.Range("J5").Value = [double]"1.488508269975597E-06"
$ws2.Range("B6").Value = [double]"3.636635927950493E-07"
$ws2.Range("C6").Value = [double]"6.31289682306857E-07"
$ws2.Range("D6").Value = 0.0130729465972379
$ws2.Range("E6").Value = 0.5840467020102458
$ws2.Range("G6").Value = 0.503298118186376
$ws2.Range("H6").Value = 0.451477748550682
$ws2.Range("I6").Value = [double]"7.069170197704366E-05"
$ws2.Range("J6").Value = [double]"1.378089287129569E-08"
$ws2.Range("B7").Value = 0.003703522545972948
$ws2.Range("C7").Value = 0.0001670144785241856
$ws2.Range("D7").Value = 0.1707566629146371
$ws2.Range("E7").Value = 0.8728760011712222
$ws2.Range("F7").Value = 0.503298118186376
$ws2.Range("H7").Value = 0.01644583639629782
$ws2.Range("I7").Value = [double]"7.649788151153913E-05"
$ws2.Range("J7").Value = 0.0001439719312719134
$ws2.Range("B8").Value = 0.01820238026848409
$ws2.Range("C8").Value = 0.0005751431850828936
$ws2.Range("D8").Value = 0.5848326098805563
$ws2.Range("E8").Value = 0.4503888688957607
$ws2.Range("F8").Value = 0.451477748550682
$ws2.Range("G8").Value = 0.01644583639629782
$ws2.Range("I8").Value = 0.0002585239466481593
$ws2.Range("J8").Value = [double]"9.806944250012606E-07"
$ws2.Range("B9").Value = 0.02441417733121787
$ws2.Range("C9").Value = 0.06538983234925433
$ws2.Range("D9").Value = 0.004547430365912764
$ws2.Range("E9").Value = 0.001294854920740995
$ws2.Range("F9").Value = [double]"7.069170197704366E-05"
$ws2.Range("G9").Value = [double]"7.649788151153913E-05"
$ws2.Range("H9").Value = 0.0002585239466481593
$ws2.Range("J9").Value = [double]"2.50769631726655E-08"
$ws2.Range("B10").Value = [double]"2.661584286300922E-10"
$ws2.Range("C10").Value = [double]"6.841640587396114E-10"
$ws2.Range("D10").Value = [double]"3.257989211746803E-09"
$ws2.Range("E10").Value = [double]"1.488508269975597E-06"
$ws2.Range("F10").Value = [double]"1.378089287129569E-08"
$ws2.Range("G10").Value = 0.0001439719312719134
$ws2.Range("H10").Value = [double]"9.806944250012606E-07"
$ws2.Range("I10").Value = [double]"2.50769631726655E-08"

# --- Sheet: Estadisticos_DM (DM test statistics) ---
$ws3 = $wb.Worksheets.Item("Estadisticos_DM")
$ws3.Range("C2").Value = 2.565571627230308
$ws3.Range("D2").Value = -3.099931797441746
$ws3.Range("E2").Value = -3.502386793273629
$ws3.Range("F2").Value = -8.947360284054245
$ws3.Range("G2").Value = -3.476598781546006
$ws3.Range("H2").Value = -2.672733267079237
$ws3.Range("I2").Value = 2.521827826554322
$ws3.Range("J2").Value = -15.75122974710012
$ws3.Range("B3").Value = -2.565571627230308
$ws3.Range("D3").Value = -3.306794187443672
$ws3.Range("E3").Value = -3.681005498625839
$ws3.Range("F3").Value = -8.544305099946431
$ws3.Range("G3").Value = -5.082431877087927
$ws3.Range("H3").Value = -4.426155284788645
$ws3.Range("I3").Value = 1.999144600078344
$ws3.Range("J3").Value = -14.66789667744244
$ws3.Range("B4").Value = 3.099931797441746
$ws3.Range("C4").Value = 3.306794187443672
$ws3.Range("E4").Value = -2.765176817227718
$ws3.Range("F4").Value = -2.841274624448338
$ws3.Range("G4").Value = -1.44395413680333
$ws3.Range("H4").Value = -0.5592419105408015
$ws3.Range("I4").Value = 3.373381154290935
$ws3.Range("J4").Value = -13.02068646215465
$ws3.Range("B5").Value = 3.502386793273629
$ws3.Range("C5").Value = 3.681005498625839
$ws3.Range("D5").Value = 2.765176817227718
$ws3.Range("F5").Value = 0.5604257539434346
$ws3.Range("G5").Value = -0.162963631956425
$ws3.Range("H5").Value = 0.7764747264221197
$ws3.Range("I5").Value = 4.008280707575426
$ws3.Range("J5").Value = -7.942848122832215
$ws3.Range("B6").Value = 8.947360284054245
$ws3.Range("C6").Value = 8.544305099946431
$ws3.Range("D6").Value = 2.841274624448338
$ws3.Range("E6").Value = -0.5604257539434346
$ws3.Range("G6").Value = -0.6870106119592574
$ws3.Range("H6").Value = 0.7745701656338183
$ws3.Range("I6").Value = 5.556784117266041
$ws3.Range("J6").Value = -11.64221546713762
$ws3.Range("B7").Value = 3.476598781546006
$ws3.Range("C7").Value = 5.082431877087927
$ws3.Range("D7").Value = 1.44395413680333
$ws3.Range("E7").Value = 0.162963631956425
$ws3.Range("F7").Value = 0.6870106119592574
$ws3.Range("H7").Value = 2.724558494569874
$ws3.Range("I7").Value = 5.512526379803321
$ws3.Range("J7").Value = -5.163178380537666
$ws3.Range("B8").Value = 2.672733267079237
$ws3.Range("C8").Value = 4.426155284788645
$ws3.Range("D8").Value = 0.5592419105408015
$ws3.Range("E8").Value = -0.7764747264221197
$ws3.Range("F8").Value = -0.7745701656338183
$ws3.Range("G8").Value = -2.724558494569874
$ws3.Range("I8").Value = 4.847387197376366
$ws3.Range("J8").Value = -8.231692865866652
$ws3.Range("B9").Value = -2.521827826554322
$ws3.Range("C9").Value = -1.999144600078344
$ws3.Range("D9").Value = -3.373381154290935
$ws3.Range("E9").Value = -4.008280707575426
$ws3.Range("F9").Value = -5.556784117266041
$ws3.Range("G9").Value = -5.512526379803321
$ws3.Range("H9").Value = -4.847387197376366
$ws3.Range("J9").Value = -11.10682590452321
$ws3.Range("B10").Value = 15.75122974710012
$ws3.Range("C10").Value = 14.66789667744244
$ws3.Range("D10").Value = 13.02068646215465
$ws3.Range("E10").Value = 7.942848122832215
$ws3.Range("F10").Value = 11.64221546713762
$ws3.Range("G10").Value = 5.163178380537666
$ws3.Range("H10").Value = 8.231692865866652
$ws3.Range("I10").Value = 11.10682590452321

# --- Sheet: Resumen (summary, re-sorted by win rate) ---
$ws4 = $wb.Worksheets.Item("Resumen")
$ws4.Range("A2").Value = "DeepAR"
$ws4.Range("B2").Value = 5
$ws4.Range("D2").Value = 3
$ws4.Range("E2").Value = 62.5
$ws4.Range("F2").Value = 0.6187332405688017
$ws4.Range("A3").Value = "Sieve Bootstrap"
$ws4.Range("B3").Value = 4
$ws4.Range("D3").Value = 4
$ws4.Range("E3").Value = 50
$ws4.Range("F3").Value = 0.6457990525229709
$ws4.Range("B4").Value = 2
$ws4.Range("D4").Value = 6
$ws4.Range("E4").Value = 25
$ws4.Range("A5").Value = "LSPMW"
$ws4.Range("B5").Value = 1
$ws4.Range("C5").Value = 1
$ws4.Range("D5").Value = 6
$ws4.Range("E5").Value = 12.5
$ws4.Range("F5").Value = 0.7349531908852894
$ws4.Range("A6").Value = "LSPM"
$ws4.Range("C6").Value = 0
$ws4.Range("D6").Value = 7
$ws4.Range("F6").Value = 0.6949129175544786
$ws4.Range("C8").Value = 2
$ws4.Range("D8").Value = 5
